$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells whose new values look numeric (e.g. "197.67", "1.00")
# as Text, so Excel keeps them as literal strings instead of silently
# converting them to floating point numbers when .Value is assigned.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.212.98'
$ws.Range("E2").Value = '  -1.47%  '

$ws.Range("D3").Value = '3.495.58'
$ws.Range("E3").Value = '  -3.33%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = '197.67'
$ws.Range("E5").Value = '  +1.54%  '

$ws.Range("D6").Value = '545.54'
$ws.Range("E6").Value = '  -6.35%  '

$ws.Range("D7").Value = '3.487.81'
$ws.Range("E7").Value = '  -3.40%  '

$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -3.12%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").Value = '0.651'
$ws.Range("E10").Value = '  -4.17%  '

$ws.Range("D11").Value = '62.59'
$ws.Range("E11").Value = '  +12.24%  '

$ws.Range("E12").Value = '  -6.92%  '

$ws.Range("E13").Value = '  -8.47%  '

$ws.Range("D14").Value = '9.73'
$ws.Range("E14").Value = '  -3.37%  '

$ws.Range("D15").Value = '4.052.44'
$ws.Range("E15").Value = '  -3.13%  '

$ws.Range("D16").Value = '3.487.59'
$ws.Range("E16").Value = '  -3.60%  '

$ws.Range("E17").Value = '  -1.78%  '

$ws.Range("D18").Value = '18.36'
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("D19").Value = '66.985.38'
$ws.Range("E19").Value = '  -1.60%  '

$ws.Range("D20").Value = '11.75'
$ws.Range("E20").Value = '  -6.24%  '

$ws.Range("E21").Value = '  -5.17%  '

$ws.Range("D22").Value = '389.68'
$ws.Range("E22").Value = '  -3.51%  '

$ws.Range("D23").Value = '3.98'
$ws.Range("E23").Value = '  -6.05%  '

$ws.Range("D24").Value = '11.79'
$ws.Range("E24").Value = '  -10.45%  '

$ws.Range("D25").Value = '81.95'
$ws.Range("E25").Value = '  -4.71%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '12.16'
$ws.Range("E26").Value = '  -3.56%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '2.78'
$ws.Range("E27").Value = '  -5.87%  '

$ws.Range("D28").Value = '3.72'
$ws.Range("E28").Value = '  -7.54%  '

$ws.Range("D29").Value = '8.72'
$ws.Range("E29").Value = '  -4.98%  '

$ws.Range("D30").Value = '30.70'

$ws.Range("D31").Value = '676.69'
$ws.Range("E31").Value = '  -1.86%  '

$ws.Range("D32").Value = '6.98'
$ws.Range("E32").Value = '  -14.67%  '

$ws.Range("D33").Value = '11.64'
$ws.Range("E33").Value = '  -4.99%  '

$ws.Range("D34").Value = '63.11'
$ws.Range("E34").Value = '  -2.61%  '

$ws.Range("E35").Value = '  -7.30%  '

$ws.Range("D36").Value = '38.56'
$ws.Range("E36").Value = '  -9.71%  '

$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("D38").Value = '0.398'
$ws.Range("E38").Value = '  -4.70%  '

$ws.Range("E39").Value = '  -3.23%  '

$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").Value = '3.051.80'
$ws.Range("E41").Value = '  -2.99%  '

$ws.Range("D42").Value = '2.97'
$ws.Range("E42").Value = '  -4.83%  '

$ws.Range("D43").Value = '0.0₃0671'
$ws.Range("E43").Value = '  -15.34%  '

$ws.Range("E44").Value = '  -13.56%  '

$ws.Range("D45").Value = '2.73'
$ws.Range("E45").Value = '  +4.43%  '

$ws.Range("D46").Value = '2.71'
$ws.Range("E46").Value = '  +5.54%  '

$ws.Range("D47").Value = '0.0395'
$ws.Range("E47").Value = '  -6.72%  '

$ws.Range("E48").Value = '  -4.44%  '

$ws.Range("D49").Value = '137.58'
$ws.Range("E49").Value = '  -4.31%  '

$ws.Range("D50").Value = '8.15'
$ws.Range("E50").Value = '  -7.93%  '

$ws.Range("D51").Value = '2.86'
$ws.Range("E51").Value = '  -8.62%  '
